$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions")

$newText = "questions = [`n    {`n        `"title`": `"Under absorption costing, your company showed a profit of `$1,200. However, the CFO of your company wants you to re-compute the profit figure using marginal costing. The closing stock contains 200 units. The absorption rate for the fixed production overheads is `$1.What is the marginal costing profit?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"`$1,000`",`n            `"`$1,200`",`n            `"`$1,400`",`n            `"`$1,600`"`n        ],`n        `"score`": `"`$1,000`"`n    },`n    {`n        `"title`": `"At the beginning of the month, your company had 40 units of material X, which were purchased at `$9/unit. The table below provides details about material purchases during the month. At the month end, 100 units were held as a closing stock. Your company applies the LIFO method for material costing.What is the cost of goods sold for the month?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"`$900`",`n            `"`$910`",`n            `"`$1,080`",`n            `"`$1,200`"`n        ],`n        `"score`": `"`$910`"`n    },`n    {`n        `"title`": `"It is the first year of operations for your company. You want to achieve at least the break-even point for the year. Your selling price is `$12/unit and the variable cost is `$8/unit. Your fixed cost for the year is estimated to be `$240,000.How many units must your company sell to achieve the break-even point?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"24,000 units`",`n            `"30,000 units`",`n            `"60,000 units`",`n            `"80,000 units`"`n        ],`n        `"score`": `"60,000 units`"`n    },`n    {`n        `"title`": `"Your company estimates a production of 400 units for the next month. Each finished unit will consume 2 lbs of material X. You have 100 lbs of material X available and you plan to maintain a closing stock of 150 lbs at the end of the month. You estimate that 50 lbs of material will be lost during the production process. The purchase price is `$2/lb.What is the budgeted material procurement cost for the next month?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"`$1,600`",`n            `"`$1,800`",`n            `"`$1,900`",`n            `"`$2,000`"`n        ],`n        `"score`": `"`$1,800`"`n    }`n]"

# Remove the second row (A2) entirely - its content moves into A1.
$ws.Rows.Item(2).Delete()

# Clear any existing formatting on A1 (bold font / border / alignment) so it
# reverts to the default style, then write the combined JSON-style text.
$cell = $ws.Range("A1")
$cell.ClearFormats()
$cell.Value = $newText
